$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the value from B1 to C1, leaving B1 empty (fix "empty column" cell issue)
$ws.Range("C1").Value2 = $ws.Range("B1").Value2
$ws.Range("B1").ClearContents()

# Update the active selection to F4 as per the recorded selection change
$ws.Range("F4").Select()
